# Tidrapport_Gr34.xlsx — "Updated hours worked on tidsrapport"
#
# Fill in the hours worked (6h each) for the four group members on the
# "v13" week sheet. The SUM formulas on v13 (B8) and on the "graf" summary
# sheet (B2/B11, and the chart built on graf!B2:B9) recalculate
# automatically. Also reflect the resulting UI state: the "v13" sheet
# becomes the active/selected tab (was "Plan"), and each sheet's selected
# cell moves to where the user was last working.

$wb = $excel.ActiveWorkbook

$plan = $wb.Worksheets.Item("Plan")
$v13  = $wb.Worksheets.Item("v13")

# Enter the hours worked by each of the 4 members for week 13.
$v13.Range("B3").Value = 6
$v13.Range("B4").Value = 6
$v13.Range("B5").Value = 6
$v13.Range("B6").Value = 6

# Plan is no longer the selected tab; its selection moved to B3.
$plan.Activate()
$plan.Range("B3").Select()

# v13 becomes the active/selected sheet, with B7 selected.
$v13.Activate()
$v13.Range("B7").Select()
